$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing table (Tabla1) and remember its reference before the edit.
$tbl = $ws.ListObjects.Item(1)

# Delete worksheet column E ("Cantidad") entirely; this shifts F->E, G->F,
# carries over the styled (underlined) empty cell from G8 to F8, and drops
# the custom column width for the removed column.
[void]$ws.Columns.Item(5).Delete()

# The ListObject/table still spans the old 5-column range (A1:E6) and keeps
# its stale "Cantidad" column definition pointing at the now-shifted data,
# so shrink it back down to the 4 real columns.
[void]$tbl.Resize($ws.Range("A1:D6"))

# Update the remaining product prices.
$ws.Range("D2").Value = 700
$ws.Range("D3").Value = 365

# Stamp a new styled (underlined), value-less marker cell at F3, matching
# the one already present at F8.
$ws.Range("F3").Font.Underline = 2

# Match the new selection recorded in the saved file.
[void]$ws.Range("G5").Select()
